$d = $word.ActiveDocument

# The "Bibliography styles" paragraph originally reads:
#   "There are various bibliography styles available. You can select the
#    style of your choice in the preamble of this document. These styles are
#    Elsevier styles based on standard styles like Harvard and Vancouver.
#    Please use BibTeX to generate your bibliography and include DOIs
#    whenever available."
#
# It needs to become:
#   "your choice in the preamble of this document. These styles are
#    Elsevier style There are various bibliography styles available. You
#    can select the style of s based on standard styles like Harvard and
#    Vancouver. Please use BibTeX to generate your bibliography and
#    include DOIs whenever available."

# --- Step 1: collapse "style of your choice in the preamble of this
#             document. These styles are Elsevier styles based on standard
#             styles like Harvard and Vancouver." down to
#             "style of s based on standard styles like Harvard and
#             Vancouver." Do this first so the later insertion (step 2)
#             cannot shift/duplicate this anchor text.
$find1 = $d.Content
$found1 = $find1.Find.Execute( `
    "style of your choice in the preamble of this document. These styles are Elsevier styles based on standard styles like Harvard and Vancouver.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "style of s based on standard styles like Harvard and Vancouver.", 1)

if (-not $found1) {
    throw "Could not find the bibliography-styles sentence to shorten."
}

# --- Step 2: locate the start of "There are various bibliography styles
#             available. You can select the" and splice new text in front
#             of it.
$find2 = $d.Content
$found2 = $find2.Find.Execute( `
    "There are various bibliography styles available. You can select the", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found2) {
    throw "Could not find the start of the bibliography-styles paragraph."
}

$insertPoint = $d.Range($find2.Start, $find2.Start)
$insertPoint.InsertBefore("your choice in the preamble of this document. These styles are Elsevier style ")
